$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TYPE")
$ws.Activate()

# Clear the value that was in A2 (keep formatting/style) and move the
# active selection to A2 (it previously pointed at C2).
$ws.Range("A2").ClearContents()
$ws.Range("A2").Select()
